$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5151125192642212
$ws.Range("B1").Value = 1.614569067955017
$ws.Range("C1").Value = 5.856376171112061
$ws.Range("D1").Value = 1.547787189483643
$ws.Range("E1").Value = 0.9460044503211975
